$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.703.29"
$ws.Range("E2").Value = "  -0.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.890.91"
$ws.Range("E3").Value = "  -1.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.54"
$ws.Range("E5").Value = "  -3.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4891"
$ws.Range("E7").Value = "  -2.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2929"
$ws.Range("E8").Value = "  -2.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06684"
$ws.Range("E9").Value = "  -2.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.888.85"
$ws.Range("E10").Value = "  -1.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.78"
$ws.Range("E11").Value = "  -1.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07244"
$ws.Range("E12").Value = "  -1.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "89.57"
$ws.Range("E13").Value = "  -2.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.025"
$ws.Range("E14").Value = "  -1.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6664"
$ws.Range("E15").Value = "  -2.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.645.84"
$ws.Range("E16").Value = "  -0.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007897"
$ws.Range("E17").Value = "  -2.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.08"
$ws.Range("E19").Value = "  -1.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.132.51"
$ws.Range("E20").Value = "  -1.05%  "

$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.751"
$ws.Range("E22").Value = "  -2.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "192.94"
$ws.Range("E23").Value = "  +5.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.090"
$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.324"
$ws.Range("E25").Value = "  -0.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.51"
$ws.Range("E26").Value = "  +2.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.33"
$ws.Range("E27").Value = "  -2.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.830"
$ws.Range("E28").Value = "  -6.13%  "

$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.275"
$ws.Range("E30").Value = "  -1.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09040"
$ws.Range("E31").Value = "  +0.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.950"
$ws.Range("E32").Value = "  -2.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05223"
$ws.Range("E33").Value = "  -0.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7350"
$ws.Range("E34").Value = "  -1.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.088"
$ws.Range("E35").Value = "  -4.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.684"
$ws.Range("E36").Value = "  +0.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01826"
$ws.Range("E37").Value = "  -6.87%  "

$ws.Range("E38").Value = "  -2.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9243"
$ws.Range("E39").Value = "  -1.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.053"
$ws.Range("E40").Value = "  -6.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4425"
$ws.Range("E41").Value = "  +0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.70"
$ws.Range("E42").Value = "  -1.73%  "

$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.744"
$ws.Range("E44").Value = "  -2.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1343"
$ws.Range("E45").Value = "  -0.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.368"
$ws.Range("E46").Value = "  -5.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4166"
$ws.Range("E47").Value = "  +5.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05829"
$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.700"
$ws.Range("E49").Value = "  +0.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.414"
$ws.Range("E50").Value = "  +1.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.28"
$ws.Range("E51").Value = "  -0.14%  "

Write-Host "Update complete"